# Final changes - 2nd June 2025
#
# Update the FlagReason sheet's column headers/messages for the Subject vs.
# Client (Buyer) round-trip flag reasons, and leave the FlagReason sheet as
# the active tab/selection (moving away from AddOpportunity).

$wb = $excel.ActiveWorkbook

$flagReason = $wb.Worksheets.Item("FlagReason")

$flagReason.Range("B2").Value = "The Subject in this engagement satisfies the requirements of a round trip. Please confirm with the deal team or CF operations team."
$flagReason.Range("B1").Value = "Subject Comment"
$flagReason.Range("C1").Value = "Client Comment"
$flagReason.Range("C2").Value = "The Buyer in this engagement satisfies the requirements of a round trip. Please confirm with the deal team or CF operations team."

[void]$flagReason.Activate()
[void]$flagReason.Range("C8").Select()
